$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(2, 7).Value = 0.4557103333333333
$ws.Cells.Item(2, 8).Value = 1.367131
$ws.Cells.Item(2, 9).Value = 0.1996401272959883
$ws.Cells.Item(2, 10).Value = 0.1996401272959883
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.05743166666666667
$ws.Cells.Item(2, 14).Value = 0.172295
$ws.Cells.Item(2, 15).Value = 0.4025715794441874
$ws.Cells.Item(2, 16).Value = 0.4025715794441875
$ws.Cells.Item(2, 17).Value = 0.02617220396055555
$ws.Cells.Item(2, 18).Value = 0.235549835645
$ws.Cells.Item(2, 19).Value = 0.08036944136598463
$ws.Cells.Item(2, 20).Value = 0.08036944136598465
$ws.Cells.Item(3, 7).Value = 0.4557103333333333
$ws.Cells.Item(3, 8).Value = 1.367131
$ws.Cells.Item(3, 9).Value = 0.1996401272959883
$ws.Cells.Item(3, 10).Value = 0.1996401272959883
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.08523033333333334
$ws.Cells.Item(3, 14).Value = 0.255691
$ws.Cells.Item(3, 15).Value = 0.5974284205558126
$ws.Cells.Item(3, 16).Value = 0.5974284205558126
$ws.Cells.Item(3, 17).Value = 0.03884034361344444
$ws.Cells.Item(3, 18).Value = 0.349563092521
$ws.Cells.Item(3, 19).Value = 0.1192706859300036
$ws.Cells.Item(3, 20).Value = 0.1192706859300036
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 9).Value = 0.2962807848215612
$ws.Cells.Item(4, 10).Value = 0.2962807848215612
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.05743166666666667
$ws.Cells.Item(4, 14).Value = 0.172295
$ws.Cells.Item(4, 15).Value = 0.4025715794441874
$ws.Cells.Item(4, 16).Value = 0.4025715794441875
$ws.Cells.Item(4, 17).Value = 0.03884149562
$ws.Cells.Item(4, 18).Value = 0.34957346058
$ws.Cells.Item(4, 19).Value = 0.1192742235045793
$ws.Cells.Item(4, 20).Value = 0.1192742235045793
$ws.Cells.Item(5, 9).Value = 0.2962807848215612
$ws.Cells.Item(5, 10).Value = 0.2962807848215612
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.08523033333333334
$ws.Cells.Item(5, 14).Value = 0.255691
$ws.Cells.Item(5, 15).Value = 0.5974284205558126
$ws.Cells.Item(5, 16).Value = 0.5974284205558126
$ws.Cells.Item(5, 17).Value = 0.057641956276
$ws.Cells.Item(5, 18).Value = 0.518777606484
$ws.Cells.Item(5, 19).Value = 0.1770065613169819
$ws.Cells.Item(5, 20).Value = 0.1770065613169819
$ws.Cells.Item(6, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.3513206666666667
$ws.Cells.Item(6, 8).Value = 1.053962
$ws.Cells.Item(6, 9).Value = 0.1539085192605057
$ws.Cells.Item(6, 10).Value = 0.1539085192605057
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.05743166666666667
$ws.Cells.Item(6, 14).Value = 0.172295
$ws.Cells.Item(6, 15).Value = 0.4025715794441874
$ws.Cells.Item(6, 16).Value = 0.4025715794441875
$ws.Cells.Item(6, 17).Value = 0.02017693142111111
$ws.Cells.Item(6, 18).Value = 0.18159238279
$ws.Cells.Item(6, 19).Value = 0.06195919568861792
$ws.Cells.Item(6, 20).Value = 0.06195919568861793
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.3513206666666667
$ws.Cells.Item(7, 8).Value = 1.053962
$ws.Cells.Item(7, 9).Value = 0.1539085192605057
$ws.Cells.Item(7, 10).Value = 0.1539085192605057
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.08523033333333334
$ws.Cells.Item(7, 14).Value = 0.255691
$ws.Cells.Item(7, 15).Value = 0.5974284205558126
$ws.Cells.Item(7, 16).Value = 0.5974284205558126
$ws.Cells.Item(7, 17).Value = 0.02994317752688889
$ws.Cells.Item(7, 18).Value = 0.269488597742
$ws.Cells.Item(7, 19).Value = 0.09194932357188779
$ws.Cells.Item(7, 20).Value = 0.09194932357188777
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.3449053333333333
$ws.Cells.Item(8, 8).Value = 1.034716
$ws.Cells.Item(8, 9).Value = 0.1510980542136751
$ws.Cells.Item(8, 10).Value = 0.1510980542136751
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.05743166666666667
$ws.Cells.Item(8, 14).Value = 0.172295
$ws.Cells.Item(8, 15).Value = 0.4025715794441874
$ws.Cells.Item(8, 16).Value = 0.4025715794441875
$ws.Cells.Item(8, 17).Value = 0.01980848813555556
$ws.Cells.Item(8, 18).Value = 0.17827639322
$ws.Cells.Item(8, 19).Value = 0.06082778233574263
$ws.Cells.Item(8, 20).Value = 0.06082778233574264
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.3449053333333333
$ws.Cells.Item(9, 8).Value = 1.034716
$ws.Cells.Item(9, 9).Value = 0.1510980542136751
$ws.Cells.Item(9, 10).Value = 0.1510980542136751
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.08523033333333334
$ws.Cells.Item(9, 14).Value = 0.255691
$ws.Cells.Item(9, 15).Value = 0.5974284205558126
$ws.Cells.Item(9, 16).Value = 0.5974284205558126
$ws.Cells.Item(9, 17).Value = 0.02939639652844445
$ws.Cells.Item(9, 18).Value = 0.264567568756
$ws.Cells.Item(9, 19).Value = 0.09027027187793245
$ws.Cells.Item(9, 20).Value = 0.09027027187793243
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 7).Value = 0.4544146666666666
$ws.Cells.Item(10, 8).Value = 1.363244
$ws.Cells.Item(10, 9).Value = 0.1990725144082698
$ws.Cells.Item(10, 10).Value = 0.1990725144082698
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.05743166666666667
$ws.Cells.Item(10, 14).Value = 0.172295
$ws.Cells.Item(10, 15).Value = 0.4025715794441874
$ws.Cells.Item(10, 16).Value = 0.4025715794441875
$ws.Cells.Item(10, 17).Value = 0.02609779166444444
$ws.Cells.Item(10, 18).Value = 0.23488012498
$ws.Cells.Item(10, 19).Value = 0.08014093654926291
$ws.Cells.Item(10, 20).Value = 0.08014093654926292
$ws.Cells.Item(11, 7).Value = 0.4544146666666666
$ws.Cells.Item(11, 8).Value = 1.363244
$ws.Cells.Item(11, 9).Value = 0.1990725144082698
$ws.Cells.Item(11, 10).Value = 0.1990725144082698
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.08523033333333334
$ws.Cells.Item(11, 14).Value = 0.255691
$ws.Cells.Item(11, 15).Value = 0.5974284205558126
$ws.Cells.Item(11, 16).Value = 0.5974284205558126
$ws.Cells.Item(11, 17).Value = 0.03872991351155555
$ws.Cells.Item(11, 18).Value = 0.348569221604
$ws.Cells.Item(11, 19).Value = 0.1189315778590069
$ws.Cells.Item(11, 20).Value = 0.1189315778590068
